$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing data row (30) into the two new rows (31, 32)
$ws.Range("A30:F30").Copy()
$ws.Range("A31:F31").PasteSpecial(-4122)
$ws.Range("A32:F32").PasteSpecial(-4122)

# Row 31: Risk Analysis
$ws.Cells.Item(31, 1).Value = "17.10.2022"
$ws.Cells.Item(31, 2).Value = 0.47916666666666669
$ws.Cells.Item(31, 3).Value = "Risk Analysis"
$ws.Cells.Item(31, 4).Value = "Documentation"
$ws.Cells.Item(31, 5).Value = 95
$ws.Cells.Item(31, 6).Value = "Collect Relevant Risk Factors, Severity, Probability, and Organise it in a Table"

# Row 32: Methodologies
$ws.Cells.Item(32, 1).Value = "17.10.2022"
$ws.Cells.Item(32, 2).Value = 0.66666666666666663
$ws.Cells.Item(32, 3).Value = "Methodologies"
$ws.Cells.Item(32, 4).Value = "Documentation"
$ws.Cells.Item(32, 5).Value = 190
$ws.Cells.Item(32, 6).Value = "Research Methdologies Suitable for Individual Projects, Document Waterfall, TDD, KanBan, and Justify Chosen Methodologies"

# Update the totals so they include the two new rows
$ws.Range("E39").Formula = "=SUM(E2:E32)"
$ws.Range("E40").Formula = "=E39 / 60"

# Update the last active selection recorded in the sheet view
[void]$ws.Range("F37").Select()
